$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style/number format/border) from B2:C2 onto B4:C4
# before writing the new values, so the "444"/"123" text values keep
# being stored as text (matching the existing style used by row 2).
$ws.Range("B2:C2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "CAJEROS"
$ws.Range("B4").Value = "444"
$ws.Range("C4").Value = "123"

$ws.Range("A3").Select()
